# Update countries & provincias Spain
# Refresh the COVID country table on sheet "Pais" with the newer snapshot
# (17 May 2020, 21:35) and re-sort a handful of countries whose total-case
# counts changed rank relative to their neighbours in the (descending by
# "Casos totales") table.
#
# Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 17 de Mayo de 2020 a las 21:35'

# --- Simple refreshes: same country keeps its row, only the counts move ---

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 1519572
$ws.Cells.Item(4, 3).Value = 11799
$ws.Cells.Item(4, 4).Value = 342658
$ws.Cells.Item(4, 5).Value = 1086190
$ws.Cells.Item(4, 7).Value = 611
$ws.Cells.Item(4, 8).Value = 90724

# Alemania (row 11)
$ws.Cells.Item(11, 2).Value = 176657
$ws.Cells.Item(11, 3).Value = 413
$ws.Cells.Item(11, 5).Value = 15220
$ws.Cells.Item(11, 7).Value = 10
$ws.Cells.Item(11, 8).Value = 8037

# Suiza (row 26)
$ws.Cells.Item(26, 4).Value = 27500
$ws.Cells.Item(26, 5).Value = 1206

# Barbados (row 170)
$ws.Cells.Item(170, 2).Value = 88
$ws.Cells.Item(170, 3).Value = 2
$ws.Cells.Item(170, 4).Value = 68
$ws.Cells.Item(170, 5).Value = 13

# --- Re-ranked block: Mayotte / El Salvador / Republica de Yibuti ---
# (Yibuti overtakes the other two, so it now leads the block in row 96)
$ws.Cells.Item(96, 1).Value = 'Republica de Yibuti'
$ws.Cells.Item(96, 2).Value = 1401
$ws.Cells.Item(96, 3).Value = 70
$ws.Cells.Item(96, 4).Value = 972
$ws.Cells.Item(96, 5).Value = 425
$ws.Cells.Item(96, 8).Value = 4

$ws.Cells.Item(97, 1).Value = 'Mayotte'
$ws.Cells.Item(97, 2).Value = 1342
$ws.Cells.Item(97, 3).Value = 30
$ws.Cells.Item(97, 4).Value = 627
$ws.Cells.Item(97, 5).Value = 697
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 18

$ws.Cells.Item(98, 1).Value = 'El Salvador'
$ws.Cells.Item(98, 2).Value = 1338
$ws.Cells.Item(98, 3).Value = 73
$ws.Cells.Item(98, 4).Value = 462
$ws.Cells.Item(98, 5).Value = 849
$ws.Cells.Item(98, 7).Value = 1
$ws.Cells.Item(98, 8).Value = 27

# --- Re-ranked block: Mali / Costa Rica ---
# (Costa Rica overtakes Mali, moves up to row 112)
$ws.Cells.Item(112, 1).Value = 'Costa Rica'
$ws.Cells.Item(112, 2).Value = 863
$ws.Cells.Item(112, 3).Value = 10
$ws.Cells.Item(112, 4).Value = 565
$ws.Cells.Item(112, 5).Value = 288
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 10

$ws.Cells.Item(113, 1).Value = 'Mali'
$ws.Cells.Item(113, 2).Value = 860
$ws.Cells.Item(113, 3).Value = 25
$ws.Cells.Item(113, 4).Value = 494
$ws.Cells.Item(113, 5).Value = 314
$ws.Cells.Item(113, 7).Value = 4
$ws.Cells.Item(113, 8).Value = 52

# --- Re-ranked block: Sudan del Sur / Ruanda ---
# (Ruanda overtakes Sudan del Sur, moves up to row 146)
$ws.Cells.Item(146, 1).Value = 'Ruanda'
$ws.Cells.Item(146, 2).Value = 292
$ws.Cells.Item(146, 3).Value = 3
$ws.Cells.Item(146, 4).Value = 197
$ws.Cells.Item(146, 5).Value = 95
$ws.Cells.Item(146, 8).Value = 0

$ws.Cells.Item(147, 1).Value = 'Sudan del Sur'
$ws.Cells.Item(147, 2).Value = 290
$ws.Cells.Item(147, 3).Value = 54
$ws.Cells.Item(147, 4).Value = 4
$ws.Cells.Item(147, 5).Value = 282
$ws.Cells.Item(147, 8).Value = 4

# --- Re-ranked block: Libia / Polinesia Francesa / Siria / Angola / Macao /
#     Zimbabue / Mauritania ---
# (Mauritania overtakes the five countries between it and Libia, so it jumps
# up to row 175 and the rest shift down one row each)
$ws.Cells.Item(175, 1).Value = 'Mauritania'
$ws.Cells.Item(175, 2).Value = 62
$ws.Cells.Item(175, 3).Value = 22
$ws.Cells.Item(175, 4).Value = 7
$ws.Cells.Item(175, 5).Value = 51
$ws.Cells.Item(175, 8).Value = 4

$ws.Cells.Item(176, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(176, 2).Value = 60
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 59
$ws.Cells.Item(176, 5).Value = 1
$ws.Cells.Item(176, 8).Value = 0

$ws.Cells.Item(177, 1).Value = 'Siria'
$ws.Cells.Item(177, 2).Value = 58
$ws.Cells.Item(177, 3).Value = 7
$ws.Cells.Item(177, 4).Value = 36
$ws.Cells.Item(177, 5).Value = 19
$ws.Cells.Item(177, 8).Value = 3

$ws.Cells.Item(178, 1).Value = 'Angola'
$ws.Cells.Item(178, 2).Value = 48
$ws.Cells.Item(178, 4).Value = 17
$ws.Cells.Item(178, 5).Value = 29
$ws.Cells.Item(178, 8).Value = 2

$ws.Cells.Item(179, 1).Value = 'Macao'
$ws.Cells.Item(179, 2).Value = 45
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 43
$ws.Cells.Item(179, 5).Value = 2
$ws.Cells.Item(179, 8).Value = 0

$ws.Cells.Item(180, 1).Value = 'Zimbabue'
$ws.Cells.Item(180, 2).Value = 44
$ws.Cells.Item(180, 3).Value = 2
$ws.Cells.Item(180, 4).Value = 17
$ws.Cells.Item(180, 5).Value = 23
